$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (changed/modified) date column for rows 2-7
# from 2023-10-13 (serial 45212) to 2023-10-22 (serial 45221)
$ws.Range("C2:C7").Value = 45221
